$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data cell D2: "no" - default formatting (written first so the shared-string
# table registers "no" before "reverses", matching the target order)
$ws.Range("D2").Value = "no"

# New header cell D1: "reverses" - bold, centered/top aligned, left+right thin border (like the
# other header cells which have a full box border, but this one only has left/right sides)
$ws.Range("D1").Value = "reverses"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("D1").VerticalAlignment = -4160     # xlTop
$ws.Range("D1").Borders.Item(7).LineStyle = 1   # xlEdgeLeft, xlContinuous
$ws.Range("D1").Borders.Item(7).Weight = 2      # xlThin
$ws.Range("D1").Borders.Item(10).LineStyle = 1  # xlEdgeRight, xlContinuous
$ws.Range("D1").Borders.Item(10).Weight = 2     # xlThin

# Move the active selection to D2, matching the saved view state
$ws.Range("D2").Select()
